$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 2694.111
$ws.Range("I2").Value = 1114.75
$ws.Range("J2").Value = 3957.6
$ws.Range("K2").Value = 1114.75
$ws.Range("L2").Value = 3957.6
$ws.Range("M2").Value = -1001.75
$ws.Range("N2").Value = -4183.6
# Row 12
$ws.Range("H12").Value = 64
$ws.Range("I12").Value = 75.666664
$ws.Range("K12").Value = 75.666664
$ws.Range("M12").Value = 94.333336
# Row 18
$ws.Range("H18").Value = 2688.2
$ws.Range("I18").Value = 2688.2
$ws.Range("K18").Value = 2688.2
$ws.Range("M18").Value = -2404.2
# Row 70
$ws.Range("H70").Value = 2000
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
# Row 73
$ws.Range("H73").Value = 2000
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
# Row 141
$ws.Range("H141").Value = 1484.5
$ws.Range("I141").Value = 1484.5
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 4453.5
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 726.5
$ws.Range("N141").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 24
$ws.Range("H24").Value = 25000
$ws.Range("J24").Value = 25000
$ws.Range("L24").Value = 25000
$ws.Range("N24").Value = -25748
# Row 32
$ws.Range("H32").Value = 7425.8423
$ws.Range("I32").Value = 6727.3335
$ws.Range("K32").Value = 6727.3335
$ws.Range("M32").Value = -6440.3335
# Row 61
$ws.Range("H61").Value = 5366.5557
$ws.Range("I61").Value = 1421.6666
$ws.Range("J61").Value = 7339
$ws.Range("K61").Value = 1421.6666
$ws.Range("L61").Value = 7339
$ws.Range("M61").Value = -1209.6666
$ws.Range("N61").Value = -7763
# Row 63
$ws.Range("H63").Value = 2000
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
# Row 66
$ws.Range("H66").Value = 2000
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
# Row 74
$ws.Range("H74").Value = 2383.3914
$ws.Range("I74").Value = 1786.4117
$ws.Range("K74").Value = 1786.4117
$ws.Range("M74").Value = -912.4117000000001
# Row 77
$ws.Range("H77").Value = 2383.3914
$ws.Range("I77").Value = 1786.4117
$ws.Range("K77").Value = 8932.058500000001
$ws.Range("M77").Value = -4564.058500000001
# Row 100
$ws.Range("H100").Value = 25000
$ws.Range("J100").Value = 25000
$ws.Range("L100").Value = 25000
$ws.Range("N100").Value = -27164
# Row 132
$ws.Range("H132").Value = 2164.647
$ws.Range("I132").Value = 2039.3572
$ws.Range("K132").Value = 6118.071599999999
$ws.Range("M132").Value = -3588.071599999999
# Row 136
$ws.Range("H136").Value = 5366.5557
$ws.Range("I136").Value = 1421.6666
$ws.Range("J136").Value = 7339
$ws.Range("K136").Value = 4264.9998
$ws.Range("L136").Value = 22017
$ws.Range("M136").Value = -1714.9998
$ws.Range("N136").Value = -27117

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 1746.3636
$ws.Range("I99").Value = 1313.5834
$ws.Range("J99").Value = 2265.7
$ws.Range("K99").Value = 1313.5834
$ws.Range("L99").Value = 2265.7
$ws.Range("M99").Value = 184.4166
$ws.Range("N99").Value = -5261.7
# Row 132
$ws.Range("H132").Value = 84999
$ws.Range("J132").Value = 84999
$ws.Range("L132").Value = 84999
$ws.Range("N132").Value = -95119
# Row 134
$ws.Range("H134").Value = 3024.7778
$ws.Range("I134").Value = 2965.375
$ws.Range("K134").Value = 8896.125
$ws.Range("M134").Value = -6361.125

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 145.27272
$ws.Range("I7").Value = 123.111115
$ws.Range("J7").Value = 245
$ws.Range("K7").Value = 123.111115
$ws.Range("L7").Value = 245
$ws.Range("M7").Value = -10.111115
$ws.Range("N7").Value = -471
# Row 22
$ws.Range("H22").Value = 28654.666
$ws.Range("I22").Value = 1179.6
$ws.Range("J22").Value = 62998.5
$ws.Range("K22").Value = 1179.6
$ws.Range("L22").Value = 62998.5
$ws.Range("M22").Value = -829.5999999999999
$ws.Range("N22").Value = -63698.5
# Row 25
$ws.Range("H25").Value = 4809.1
$ws.Range("I25").Value = 2000
$ws.Range("K25").Value = 2000
$ws.Range("M25").Value = -1826
# Row 51
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
# Row 59
$ws.Range("H59").Value = 112271.43
$ws.Range("I59").Value = 100000
$ws.Range("J59").Value = 114316.664
$ws.Range("K59").Value = 100000
$ws.Range("L59").Value = 114316.664
$ws.Range("M59").Value = -98855
$ws.Range("N59").Value = -116606.664
# Row 60
$ws.Range("H60").Value = 15805.25
$ws.Range("I60").Value = 9407
$ws.Range("J60").Value = 35000
$ws.Range("K60").Value = 9407
$ws.Range("L60").Value = 35000
$ws.Range("M60").Value = -8896
$ws.Range("N60").Value = -36022
# Row 61
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
# Row 104
$ws.Range("H104").Value = 50000
$ws.Range("J104").Value = 50000
$ws.Range("L104").Value = 50000
$ws.Range("N104").Value = -55242
# Row 122
$ws.Range("H122").Value = 1592.4286
$ws.Range("I122").Value = 1222.2858
$ws.Range("J122").Value = 1962.5714
$ws.Range("K122").Value = 3666.8574
$ws.Range("L122").Value = 5887.7142
$ws.Range("M122").Value = -1216.8574
$ws.Range("N122").Value = -10787.7142
# Row 132
$ws.Range("H132").Value = 2037.1666
$ws.Range("I132").Value = 1555.9231
$ws.Range("K132").Value = 4667.7693
$ws.Range("M132").Value = -2137.7693

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 112
$ws.Range("H112").Value = 7641
$ws.Range("I112").Value = 5189.6665
$ws.Range("J112").Value = 14995
$ws.Range("K112").Value = 15568.9995
$ws.Range("L112").Value = 44985
$ws.Range("M112").Value = -14460.9995
$ws.Range("N112").Value = -47201
# Row 132
$ws.Range("H132").Value = 3763.9167
$ws.Range("J132").Value = 2512
$ws.Range("L132").Value = 22608
$ws.Range("N132").Value = -27668
# Row 133
$ws.Range("H133").Value = 3995
$ws.Range("I133").Value = 3995
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 11985
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -6925
$ws.Range("N133").ClearContents()
# Row 134
$ws.Range("H134").Value = 12231.363
$ws.Range("I134").Value = 1573.5
$ws.Range("J134").Value = 18321.572
$ws.Range("K134").Value = 4720.5
$ws.Range("L134").Value = 54964.716
$ws.Range("M134").Value = 349.5
$ws.Range("N134").Value = -65104.716
# Row 138
$ws.Range("H138").Value = 3199.8
$ws.Range("I138").Value = 3199.8
$ws.Range("K138").Value = 9599.400000000001
$ws.Range("M138").Value = -4459.400000000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 1080
$ws.Range("I132").Value = 1000
$ws.Range("K132").Value = 3000
$ws.Range("M132").Value = -470

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 166666.67
$ws.Range("I2").Value = 166666.67
$ws.Range("K2").Value = 166666.67
$ws.Range("M2").Value = -166554.67
# Row 7
$ws.Range("H7").Value = 7659.391
$ws.Range("I7").Value = 3813.3333
$ws.Range("K7").Value = 3813.3333
$ws.Range("M7").Value = -3701.3333
# Row 93
$ws.Range("H93").Value = 2016.9412
$ws.Range("I93").Value = 1754.4546
$ws.Range("K93").Value = 1754.4546
$ws.Range("M93").Value = -506.4546
# Row 116
$ws.Range("H116").Value = 150000
$ws.Range("J116").Value = 150000
$ws.Range("L116").Value = 150000
$ws.Range("N116").Value = -159178
# Row 122
$ws.Range("H122").Value = 6195
$ws.Range("I122").Value = 5749.4287
$ws.Range("J122").Value = 6888.1113
$ws.Range("K122").Value = 17248.2861
$ws.Range("L122").Value = 20664.3339
$ws.Range("M122").Value = -14798.2861
$ws.Range("N122").Value = -25564.3339
# Row 126
$ws.Range("H126").Value = 7659.391
$ws.Range("I126").Value = 3813.3333
$ws.Range("K126").Value = 11439.9999
$ws.Range("M126").Value = -8969.999899999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 1125000
$ws.Range("J2").Value = 250000
$ws.Range("L2").Value = 250000
$ws.Range("N2").Value = -250224
# Row 75
$ws.Range("H75").Value = 87118
$ws.Range("I75").Value = 87118
$ws.Range("K75").Value = 87118
$ws.Range("M75").Value = -86182
# Row 78
$ws.Range("H78").Value = 87118
$ws.Range("I78").Value = 87118
$ws.Range("K78").Value = 261354
$ws.Range("M78").Value = -256674
# Row 98
$ws.Range("H98").Value = 32333
$ws.Range("J98").Value = 32333
$ws.Range("L98").Value = 32333
$ws.Range("N98").Value = -38323
# Row 122
$ws.Range("H122").Value = 652
$ws.Range("I122").Value = 652
$ws.Range("K122").Value = 1956
$ws.Range("M122").Value = 494
# Row 126
$ws.Range("H126").Value = 8749.333000000001
$ws.Range("I126").Value = 8000
$ws.Range("K126").Value = 24000
$ws.Range("M126").Value = -21530
